# Commit: "git commit calificacion hasta p23"
# Adds a new repeat-group column "l1" to the p26 repeat-group export sheet.
#
# Net effect observed in the target workbook:
#   - A new shared string "l1" is introduced.
#   - Header cell AC1 (previously "nota_iniciativa") now reads "l1".
#   - A brand-new header cell AD1 is added, reading "nota_iniciativa"
#     (i.e. the old header text moved one column to the right).
#   - Every data row (2-97) gets a new AD cell holding the numeric value 0,
#     matching the existing AC/AB/AA/Z "0" flag columns.
#   - The used range grows from A1:AC97 to A1:AD97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
# AC1 flips from "nota_iniciativa" to "l1"...
$ws.Range("AC1").Value = "l1"

# ...and the new AD1 header takes over the old "nota_iniciativa" text.
$ws.Range("AD1").Value = "nota_iniciativa"

# Match AD1's look (bold, centered, bordered) to the rest of the header row.
$headerSample = $ws.Range("AC1")
$newHeader = $ws.Range("AD1")
$newHeader.Font.Bold = $headerSample.Font.Bold
$newHeader.HorizontalAlignment = $headerSample.HorizontalAlignment
$newHeader.VerticalAlignment = $headerSample.VerticalAlignment
$newHeader.Borders.LineStyle = $headerSample.Borders.LineStyle

# --- Data rows --------------------------------------------------------------
# Rows 2-97 each get a new AD cell set to 0, same as the neighbouring
# Z/AA/AB/AC "0" marker columns.
$lastRow = 97
$ws.Range("AD2:AD" + $lastRow).Value = 0
